$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2830.7273
$ws.Range("I2").Value = 924.8333
$ws.Range("K2").Value = 924.8333
$ws.Range("M2").Value = -811.8333

$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

$ws.Range("H32").Value = 1045.3077

$ws.Range("H42").Value = 1506
$ws.Range("I42").Value = 25.666666
$ws.Range("J42").Value = 2986.3333
$ws.Range("K42").Value = 76.99999800000001
$ws.Range("L42").Value = 8958.999899999999
$ws.Range("M42").Value = 153.000002
$ws.Range("N42").Value = -9418.999899999999

$ws.Range("H58").Value = 1745.625
$ws.Range("I58").Value = 991.3333
$ws.Range("J58").Value = 4008.5
$ws.Range("K58").Value = 2973.9999
$ws.Range("L58").Value = 12025.5
$ws.Range("M58").Value = -2823.9999
$ws.Range("N58").Value = -12325.5

$ws.Range("H93").Value = 19515.8
$ws.Range("J93").Value = 19515.8
$ws.Range("L93").Value = 19515.8
$ws.Range("N93").Value = -24507.8

$ws.Range("H135").Value = 704.2
$ws.Range("J135").Value = 1099.5
$ws.Range("L135").Value = 9895.5
$ws.Range("N135").Value = -14965.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2303.35
$ws.Range("I32").Value = 2124.1553
$ws.Range("K32").Value = 2124.1553
$ws.Range("M32").Value = -1837.1553

$ws.Range("H45").Value = 1857.6666
$ws.Range("I45").Value = 1677.2941
$ws.Range("J45").Value = 2624.25
$ws.Range("K45").Value = 1677.2941
$ws.Range("L45").Value = 2624.25
$ws.Range("M45").Value = -1300.2941
$ws.Range("N45").Value = -3378.25

$ws.Range("H122").Value = 1644.375
$ws.Range("I122").Value = 1367.5
$ws.Range("J122").Value = 2475
$ws.Range("K122").Value = 4102.5
$ws.Range("L122").Value = 7425
$ws.Range("M122").Value = -1652.5
$ws.Range("N122").Value = -12325

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H105").Value = 4890.8335
$ws.Range("I105").Value = 5486.75
$ws.Range("J105").Value = 3699
$ws.Range("K105").Value = 5486.75
$ws.Range("L105").Value = 3699
$ws.Range("M105").Value = -3739.75
$ws.Range("N105").Value = -7193

$ws.Range("H107").Value = 1302.5714
$ws.Range("I107").Value = 923.6
$ws.Range("K107").Value = 923.6
$ws.Range("M107").Value = 996.4

$ws.Range("H134").Value = 8352.9375
$ws.Range("I134").Value = 8643.134
$ws.Range("K134").Value = 25929.402
$ws.Range("M134").Value = -23394.402

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 562.25
$ws.Range("I105").Value = 374.5
$ws.Range("J105").Value = 750
$ws.Range("K105").Value = 374.5
$ws.Range("L105").Value = 750
$ws.Range("M105").Value = 1372.5
$ws.Range("N105").Value = -4244

$ws.Range("H134").Value = 1998
$ws.Range("I134").Value = 1998
$ws.Range("K134").Value = 5994
$ws.Range("M134").Value = -3459

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3038.2666
$ws.Range("J34").Value = 3169.5715
$ws.Range("L34").Value = 9508.7145
$ws.Range("N34").Value = -9676.7145

$ws.Range("H80").Value = 12981.866
$ws.Range("J80").Value = 13456.077
$ws.Range("L80").Value = 40368.231
$ws.Range("N80").Value = -42240.231

$ws.Range("H83").Value = 12981.866
$ws.Range("J83").Value = 13456.077
$ws.Range("L83").Value = 121104.693
$ws.Range("N83").Value = -130464.693

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 289.8
$ws.Range("I2").Value = 283.1111
$ws.Range("K2").Value = 283.1111
$ws.Range("M2").Value = -170.1111

$ws.Range("H18").Value = 37000
$ws.Range("J18").Value = 37000
$ws.Range("L18").Value = 37000
$ws.Range("N18").Value = -37586

$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H57").Value = 23132.715
$ws.Range("J57").Value = 24988.166
$ws.Range("L57").Value = 24988.166
$ws.Range("N57").Value = -26628.166

$ws.Range("H70").Value = 37042424
$ws.Range("J70").Value = 10000
$ws.Range("L70").Value = 10000
$ws.Range("N70").Value = -10540

$ws.Range("H73").Value = 37042424
$ws.Range("J73").Value = 10000
$ws.Range("L73").Value = 10000
$ws.Range("N73").Value = -11872

$ws.Range("H80").Value = 1798.625
$ws.Range("J80").Value = 2472
$ws.Range("L80").Value = 2472
$ws.Range("N80").Value = -4468

$ws.Range("H83").Value = 1798.625
$ws.Range("J83").Value = 2472
$ws.Range("L83").Value = 12360
$ws.Range("N83").Value = -22344

$ws.Range("H113").Value = 738
$ws.Range("I113").Value = 738
$ws.Range("K113").Value = 738
$ws.Range("M113").Value = 1432

$ws.Range("H126").Value = 9492.286
$ws.Range("I126").Value = 9949.299999999999
$ws.Range("J126").Value = 8349.75
$ws.Range("K126").Value = 29847.9
$ws.Range("L126").Value = 25049.25
$ws.Range("M126").Value = -27377.9
$ws.Range("N126").Value = -29989.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 992.3333
$ws.Range("J22").Value = 1148.5
$ws.Range("L22").Value = 1148.5
$ws.Range("N22").Value = -1738.5

$ws.Range("H27").Value = 992.3333
$ws.Range("J27").Value = 1148.5
$ws.Range("L27").Value = 1148.5
$ws.Range("N27").Value = -1362.5

$ws.Range("H46").Value = 2999.9656
$ws.Range("J46").Value = 3923.077
$ws.Range("L46").Value = 3923.077
$ws.Range("N46").Value = -4299.077

$ws.Range("H55").Value = 166.42308
$ws.Range("I55").Value = 137.5625
$ws.Range("J55").Value = 212.6
$ws.Range("K55").Value = 137.5625
$ws.Range("L55").Value = 212.6
$ws.Range("M55").Value = 35.4375
$ws.Range("N55").Value = -558.6

$ws.Range("H61").Value = 8251.25
$ws.Range("I61").Value = 7750.5
$ws.Range("J61").Value = 8752
$ws.Range("K61").Value = 7750.5
$ws.Range("L61").Value = 8752
$ws.Range("M61").Value = -7548.5
$ws.Range("N61").Value = -9156

$ws.Range("H113").Value = 8251.25
$ws.Range("I113").Value = 7750.5
$ws.Range("J113").Value = 8752
$ws.Range("K113").Value = 7750.5
$ws.Range("L113").Value = 8752
$ws.Range("M113").Value = -5580.5
$ws.Range("N113").Value = -13092

$ws.Range("H122").Value = 3504
$ws.Range("I122").Value = 3504
$ws.Range("K122").Value = 10512
$ws.Range("M122").Value = -8062

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2761.7856
$ws.Range("J81").Value = 1945
$ws.Range("L81").Value = 3890
$ws.Range("N81").Value = -6012

$ws.Range("H84").Value = 2761.7856
$ws.Range("J84").Value = 1945
$ws.Range("L84").Value = 19450
$ws.Range("N84").Value = -30058

$ws.Range("H95").Value = 28133.334
$ws.Range("J95").Value = 28133.334
$ws.Range("L95").Value = 28133.334
$ws.Range("N95").Value = -33625.334

$ws.Range("H122").Value = 2491.4666
$ws.Range("I122").Value = 2291.0833
$ws.Range("K122").Value = 6873.249899999999
$ws.Range("M122").Value = -4423.249899999999

$ws.Range("H136").Value = 3653.0488
$ws.Range("I136").Value = 3520.8108
$ws.Range("J136").Value = 4876.25
$ws.Range("K136").Value = 10562.4324
$ws.Range("L136").Value = 14628.75
$ws.Range("M136").Value = -8012.432400000002
$ws.Range("N136").Value = -19728.75
